$p = $ppt.ActivePresentation

# Walk every slide and every shape looking for tables, then push the
# writer's internally-computed column widths onto them (in points;
# 2514600 EMU == 198 pt) instead of the old evenly-split values.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $tbl.Columns.Item($c).Width = 198
            }
        }
    }
}
